$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<do>"
$ws.Range("C2").Value = 31

# Row 3
$ws.Range("B3").Value = "<conta>"
$ws.Range("C3").Value = 30

# Row 4
$ws.Range("B4").Value = "<down>"
$ws.Range("C4").Value = 29

# Row 5
$ws.Range("B5").Value = "<quet>"
$ws.Range("C5").Value = 29

# Row 6
$ws.Range("B6").Value = "<now>"
$ws.Range("C6").Value = 27

# Row 7
$ws.Range("C7").Value = 28

# Row 8
$ws.Range("C8").Value = 26

# Row 10
$ws.Range("C10").Value = 24

# Row 11
$ws.Range("B11").Value = "<enter>"
$ws.Range("C11").Value = 24

# Row 12
$ws.Range("C12").Value = 30

# Row 13
$ws.Range("C13").Value = 35

# Row 14
$ws.Range("B14").Value = "<sie>"
$ws.Range("C14").Value = 34

# Row 15
$ws.Range("B15").Value = "<in>"
$ws.Range("C15").Value = 24

# Row 16
$ws.Range("C16").Value = 32

# Row 17
$ws.Range("B17").Value = "<so>"
$ws.Range("C17").Value = 26

# Row 18
$ws.Range("C18").Value = 29
